$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("C4").Value = "北京·原神庙会 .  海灯节同人展（取消）"
$ws.Range("F4").Value = 1291
$ws.Range("F6").Value = 316
$ws.Range("F7").Value = 1120
$ws.Range("F9").Value = 6965
$ws.Range("F13").Value = 7861
$ws.Range("F16").Value = 5454
$ws.Range("F17").Value = 45
$ws.Range("F18").Value = 2330
$ws.Range("F19").Value = 984
$ws.Range("F21").Value = 276
$ws.Range("F22").Value = 371
$ws.Range("F25").Value = 324
$ws.Range("F28").Value = 2120
$ws.Range("F31").Value = 69
$ws.Range("F32").Value = 64
$ws.Range("F35").Value = 23
$ws.Range("F36").Value = 1426
$ws.Range("F39").Value = 2165
$ws.Range("F40").Value = 2185
$ws.Range("G4").Value = "不可售"

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 65
$ws.Range("F4").Value = 39

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 251
$ws.Range("F3").Value = 1264

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("C6").Value = "北京·原神庙会 .  海灯节同人展（取消）"
$ws.Range("F2").Value = 251
$ws.Range("F4").Value = 1264
$ws.Range("F6").Value = 1291
$ws.Range("F9").Value = 316
$ws.Range("F10").Value = 1120
$ws.Range("F12").Value = 6965
$ws.Range("F16").Value = 7861
$ws.Range("F19").Value = 5454
$ws.Range("F20").Value = 45
$ws.Range("F21").Value = 2330
$ws.Range("F22").Value = 984
$ws.Range("F24").Value = 276
$ws.Range("F25").Value = 371
$ws.Range("F27").Value = 65
$ws.Range("F29").Value = 39
$ws.Range("F30").Value = 324
$ws.Range("F33").Value = 2120
$ws.Range("F36").Value = 69
$ws.Range("F37").Value = 64
$ws.Range("F40").Value = 23
$ws.Range("F42").Value = 1426
$ws.Range("F45").Value = 2165
$ws.Range("F47").Value = 2185
$ws.Range("G6").Value = "不可售"
